{"js": "// The document stores pseudo-XML markup as literal text runs, e.g.\n// \"<id>p089v_a1</id>\" split across three runs (opening tag run with\n// Courier New formatting, a plain-formatted id-value run, and a closing\n// tag run with Courier New formatting). The edit collapses each triple\n// into a single run (keeping the first run's formatting) whose text is\n// the updated id, \"p089v_1\" / \"p089v_2\" (dropping the \"a\").\nconst replacements = [\n  [\"<id>p089v_a1</id>\", \"<id>p089v_1</id>\"],\n  [\"<id>p089v_a2</id>\", \"<id>p089v_2</id>\"],\n];\n\nfor (const [needle, replacement] of replacements) {\n  const results = context.document.body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The document stores pseudo-XML markup as literal text runs, e.g.\n# \"<id>p089v_a1</id>\" split across three runs (opening tag run with\n# Courier New formatting, a plain-formatted id-value run, and a closing\n# tag run with Courier New formatting). The edit collapses each triple\n# into a single run (keeping the first run's formatting) whose text is\n# the updated id, \"p089v_1\" / \"p089v_2\" (dropping the \"a\").\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"<id>p089v_a1</id>\", \"<id>p089v_1</id>\"),\n    @(\"<id>p089v_a2</id>\", \"<id>p089v_2</id>\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
